# Updates cryptos list values (price + volume) per upstream scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.861.90"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.223.57"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.25%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.90"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.34%  "
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.598"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +12.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0968"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.27"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.20"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.97%  "
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "2.554.43"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.866"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").Value = "2.221.26"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "41.722.16"
$ws.Range("E19").Value = "  -1.33%  "
$ws.Range("D20").Value = "0.0₃0967"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("E21").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.43%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +6.00%  "
$ws.Range("E28").Value = "  +7.09%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.19"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.81%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.22"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.75%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("E32").Value = "  +2.99%  "
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.125"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.32%  "
$ws.Range("E35").Value = "  +0.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.73"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "25.60"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +16.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.02"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.37%  "
$ws.Range("E39").Value = "  +7.52%  "
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "66.74"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.83%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +19.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.203"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  -4.50%  "
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("E50").Value = "  +6.44%  "
$ws.Range("E51").Value = "  +1.16%  "

Write-Output "Applied all changes"
